$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1), columns BA:BU ---
$ws.Range("BA1").Value = "APP"
$ws.Range("BB1").Value = "UX_A"
$ws.Range("BC1").Value = "UX_P"
$ws.Range("BD1").Value = "UX_E"
$ws.Range("BE1").Value = "UX_D"
$ws.Range("BF1").Value = "UX_S"
$ws.Range("BG1").Value = "UX_N"
$ws.Range("BH1").Value = "SUS"
$ws.Range("BI1").Value = "CL_W_MD"
$ws.Range("BJ1").Value = "CL_W_PD"
$ws.Range("BK1").Value = "CL_W_TD"
$ws.Range("BL1").Value = "CL_W_E"
$ws.Range("BM1").Value = "CL_W_P"
$ws.Range("BN1").Value = "CL_W_F"
$ws.Range("BO1").Value = "CL_MD"
$ws.Range("BP1").Value = "CL_PD"
$ws.Range("BQ1").Value = "CL_TD"
$ws.Range("BR1").Value = "CL_E"
$ws.Range("BS1").Value = "CL_P"
$ws.Range("BT1").Value = "CL_F"
$ws.Range("BU1").Value = "CL_SCORE"

# --- Update data rows 2-9, columns BA:BU with new FORM results ---
# Row 2
$ws.Range("BA2").Value = 0
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BD2").Value = 0.25
$ws.Range("BE2").Value = -0.75
$ws.Range("BF2").Value = 0.25
$ws.Range("BG2").Value = 1.5
$ws.Range("BH2").Value = 90
$ws.Range("BI2").Value = 0.133
$ws.Range("BJ2").Value = 0.2
$ws.Range("BK2").Value = 0.067
$ws.Range("BL2").Value = 0.2
$ws.Range("BM2").Value = 0.067
$ws.Range("BN2").Value = 0.333
$ws.Range("BO2").Value = 0.4
$ws.Range("BP2").Value = 1.4
$ws.Range("BQ2").Value = 0.133
$ws.Range("BR2").Value = 0.8
$ws.Range("BS2").Value = 0.6
$ws.Range("BT2").Value = 0.333
$ws.Range("BU2").Value = 4.666

# Row 3
$ws.Range("BA3").Value = 0
$ws.Range("BB3").Value = 0
$ws.Range("BC3").Value = 0
$ws.Range("BD3").Value = 0.25
$ws.Range("BE3").Value = -0.75
$ws.Range("BF3").Value = 0.25
$ws.Range("BG3").Value = 1.5
$ws.Range("BH3").Value = 90
$ws.Range("BI3").Value = 0.133
$ws.Range("BJ3").Value = 0.2
$ws.Range("BK3").Value = 0.067
$ws.Range("BL3").Value = 0.2
$ws.Range("BM3").Value = 0.067
$ws.Range("BN3").Value = 0.333
$ws.Range("BO3").Value = 0.4
$ws.Range("BP3").Value = 1.4
$ws.Range("BQ3").Value = 0.133
$ws.Range("BR3").Value = 0.8
$ws.Range("BS3").Value = 0.6
$ws.Range("BT3").Value = 0.333
$ws.Range("BU3").Value = 4.666

# Row 4
$ws.Range("BA4").Value = 0
$ws.Range("BB4").Value = 0
$ws.Range("BC4").Value = 0
$ws.Range("BD4").Value = 0.25
$ws.Range("BE4").Value = -0.75
$ws.Range("BF4").Value = 0.25
$ws.Range("BG4").Value = 1.5
$ws.Range("BH4").Value = 90
$ws.Range("BI4").Value = 0.133
$ws.Range("BJ4").Value = 0.2
$ws.Range("BK4").Value = 0.067
$ws.Range("BL4").Value = 0.2
$ws.Range("BM4").Value = 0.067
$ws.Range("BN4").Value = 0.333
$ws.Range("BO4").Value = 0.4
$ws.Range("BP4").Value = 1.4
$ws.Range("BQ4").Value = 0.133
$ws.Range("BR4").Value = 0.8
$ws.Range("BS4").Value = 0.6
$ws.Range("BT4").Value = 0.333
$ws.Range("BU4").Value = 4.666

# Row 5
$ws.Range("BA5").Value = 0
$ws.Range("BB5").Value = 0
$ws.Range("BC5").Value = 0
$ws.Range("BD5").Value = 0.25
$ws.Range("BE5").Value = -0.75
$ws.Range("BF5").Value = 0.25
$ws.Range("BG5").Value = 1.5
$ws.Range("BH5").Value = 90
$ws.Range("BI5").Value = 0.133
$ws.Range("BJ5").Value = 0.2
$ws.Range("BK5").Value = 0.067
$ws.Range("BL5").Value = 0.2
$ws.Range("BM5").Value = 0.067
$ws.Range("BN5").Value = 0.333
$ws.Range("BO5").Value = 0.4
$ws.Range("BP5").Value = 1.4
$ws.Range("BQ5").Value = 0.133
$ws.Range("BR5").Value = 0.8
$ws.Range("BS5").Value = 0.6
$ws.Range("BT5").Value = 0.333
$ws.Range("BU5").Value = 4.666

# Row 6
$ws.Range("BA6").Value = 0
$ws.Range("BB6").Value = -0.333
$ws.Range("BC6").Value = 0
$ws.Range("BD6").Value = 0.25
$ws.Range("BE6").Value = -0.25
$ws.Range("BF6").Value = -1
$ws.Range("BG6").Value = -0.25
$ws.Range("BH6").Value = 92.5
$ws.Range("BI6").Value = 0.333
$ws.Range("BJ6").Value = 0.067
$ws.Range("BK6").Value = 0.133
$ws.Range("BL6").Value = 0.2
$ws.Range("BM6").Value = 0.267
$ws.Range("BN6").Value = 0
$ws.Range("BO6").Value = 0.333
$ws.Range("BP6").Value = 0.067
$ws.Range("BQ6").Value = 0.133
$ws.Range("BR6").Value = 0.4
$ws.Range("BS6").Value = 2.133
$ws.Range("BT6").Value = 0
$ws.Range("BU6").Value = 4.066

# Row 7
$ws.Range("BA7").Value = 0
$ws.Range("BB7").Value = -0.333
$ws.Range("BC7").Value = 0
$ws.Range("BD7").Value = 0.25
$ws.Range("BE7").Value = -0.25
$ws.Range("BF7").Value = -1
$ws.Range("BG7").Value = -0.25
$ws.Range("BH7").Value = 92.5
$ws.Range("BI7").Value = 0.333
$ws.Range("BJ7").Value = 0.067
$ws.Range("BK7").Value = 0.133
$ws.Range("BL7").Value = 0.2
$ws.Range("BM7").Value = 0.267
$ws.Range("BN7").Value = 0
$ws.Range("BO7").Value = 0.333
$ws.Range("BP7").Value = 0.067
$ws.Range("BQ7").Value = 0.133
$ws.Range("BR7").Value = 0.4
$ws.Range("BS7").Value = 2.133
$ws.Range("BT7").Value = 0
$ws.Range("BU7").Value = 4.066

# Row 8
$ws.Range("BA8").Value = 0
$ws.Range("BB8").Value = -0.333
$ws.Range("BC8").Value = 0
$ws.Range("BD8").Value = 0.25
$ws.Range("BE8").Value = -0.25
$ws.Range("BF8").Value = -1
$ws.Range("BG8").Value = -0.25
$ws.Range("BH8").Value = 92.5
$ws.Range("BI8").Value = 0.333
$ws.Range("BJ8").Value = 0.067
$ws.Range("BK8").Value = 0.133
$ws.Range("BL8").Value = 0.2
$ws.Range("BM8").Value = 0.267
$ws.Range("BN8").Value = 0
$ws.Range("BO8").Value = 0.333
$ws.Range("BP8").Value = 0.067
$ws.Range("BQ8").Value = 0.133
$ws.Range("BR8").Value = 0.4
$ws.Range("BS8").Value = 2.133
$ws.Range("BT8").Value = 0
$ws.Range("BU8").Value = 4.066

# Row 9
$ws.Range("BA9").Value = 0
$ws.Range("BB9").Value = -0.333
$ws.Range("BC9").Value = 0
$ws.Range("BD9").Value = 0.25
$ws.Range("BE9").Value = -0.25
$ws.Range("BF9").Value = -1
$ws.Range("BG9").Value = -0.25
$ws.Range("BH9").Value = 92.5
$ws.Range("BI9").Value = 0.333
$ws.Range("BJ9").Value = 0.067
$ws.Range("BK9").Value = 0.133
$ws.Range("BL9").Value = 0.2
$ws.Range("BM9").Value = 0.267
$ws.Range("BN9").Value = 0
$ws.Range("BO9").Value = 0.333
$ws.Range("BP9").Value = 0.067
$ws.Range("BQ9").Value = 0.133
$ws.Range("BR9").Value = 0.4
$ws.Range("BS9").Value = 2.133
$ws.Range("BT9").Value = 0
$ws.Range("BU9").Value = 4.066

# --- Remove now-unused trailing columns BV:CW ---
$ws.Columns("BV:CW").Delete()
